$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - new calibration values (Pre-Cond - 2)
$ws.Range("B4").Value = 26459.3715
$ws.Range("C4").Value = 119.01181
$ws.Range("D4").Value = -1.06083
$ws.Range("E4").Value = 0.00247
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 7901.94755

# Row 5 - new calibration values (Post-Cond - 2)
$ws.Range("B5").Value = 21939.00353
$ws.Range("C5").Value = 245.64123
$ws.Range("D5").Value = -0.99305
$ws.Range("E5").Value = 0.0061
$ws.Range("F5").Value = 0.96469
$ws.Range("G5").Value = 42048.43112

# Match the existing cell style in these rows (center/center alignment, no bold)
$ws.Range("B4:G5").HorizontalAlignment = -4108
$ws.Range("B4:G5").VerticalAlignment = -4108

# Update the active selection to match the saved view state
$ws.Range("F10").Select()
